# Apply updated cryptocurrency price/volume data per commit
# "Updated cryptos list on Sat Sep 16 23:04:31 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'26.771.92"
$ws.Range('E2').Value = "'  -0.07%  "
# Row 3
$ws.Range('D3').Value = "'1.648.78"
$ws.Range('E3').Value = "'  -0.11%  "
# Row 4
$ws.Range('E4').Value = "'  +0.83%  "
# Row 5
$ws.Range('D5').Value = "'216.70"
$ws.Range('E5').Value = "'  +0.90%  "
# Row 6
$ws.Range('E6').Value = "'  +0.24%  "
# Row 7
$ws.Range('E7').Value = "'  +0.73%  "
# Row 8
$ws.Range('E8').Value = "'  +0.02%  "
# Row 9
$ws.Range('D9').Value = "'0.0626"
$ws.Range('E9').Value = "'  -0.23%  "
# Row 10
$ws.Range('E10').Value = "'  +0.08%  "
# Row 11
$ws.Range('E11').Value = "'  +0.17%  "
# Row 12
$ws.Range('D12').Value = "'1.873.50"
$ws.Range('E12').Value = "'  -0.21%  "
# Row 13
$ws.Range('B13').Value = "'WrappedEther"
$ws.Range('C13').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D13').Value = "'1.652.09"
$ws.Range('E13').Value = "'  +0.05%  "
# Row 14
$ws.Range('B14').Value = "'Polkadot"
$ws.Range('C14').Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('D14').Value = "'4.21"
$ws.Range('E14').Value = "'  +0.97%  "
# Row 15
$ws.Range('E15').Value = "'  -0.10%  "
# Row 16
$ws.Range('D16').Value = "'65.56"
$ws.Range('E16').Value = "'  -0.58%  "
# Row 17
$ws.Range('D17').Value = "'26.789.90"
$ws.Range('E17').Value = "'  +0.04%  "
# Row 18
$ws.Range('D18').Value = "'0.0₃0744"
$ws.Range('E18').Value = "'  -0.34%  "
# Row 19
$ws.Range('D19').Value = "'217.35"
$ws.Range('E19').Value = "'  -0.74%  "
# Row 20
$ws.Range('E20').Value = "'  +0.77%  "
# Row 21
$ws.Range('B21').Value = "'Uniswap"
$ws.Range('C21').Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range('D21').Value = "'4.38"
$ws.Range('E21').Value = "'  +0.41%  "
# Row 22
$ws.Range('B22').Value = "'Toncoin"
$ws.Range('C22').Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('D22').Value = "'2.50"
$ws.Range('E22').Value = "'  +18.38%  "
# Row 23
$ws.Range('E23').Value = "'  -0.66%  "
# Row 24
$ws.Range('D24').Value = "'9.46"
$ws.Range('E24').Value = "'  -0.17%  "
# Row 25
$ws.Range('D25').Value = "'145.63"
$ws.Range('E25').Value = "'  -1.44%  "
# Row 26
$ws.Range('E26').Value = "'  +0.80%  "
# Row 27
$ws.Range('D27').Value = "'0.120"
$ws.Range('E27').Value = "'  -0.46%  "
# Row 28
$ws.Range('E28').Value = "'  +3.67%  "
# Row 29
$ws.Range('E29').Value = "'  -0.35%  "
# Row 30
$ws.Range('D30').Value = "'0.0522"
$ws.Range('E30').Value = "'  +0.44%  "
# Row 31
$ws.Range('E31').Value = "'  +0.89%  "
# Row 32
$ws.Range('D32').Value = "'3.35"
$ws.Range('E32').Value = "'  -0.91%  "
# Row 33
$ws.Range('E33').Value = "'  -0.41%  "
# Row 34
$ws.Range('E34').Value = "'  +1.35%  "
# Row 35
$ws.Range('D35').Value = "'1.277.18"
$ws.Range('E35').Value = "'  +0.37%  "
# Row 36
$ws.Range('E36').Value = "'  +2.08%  "
# Row 37
$ws.Range('E37').Value = "'  +1.00%  "
# Row 38
$ws.Range('E38').Value = "'  +5.32%  "
# Row 39
$ws.Range('E39').Value = "'  +2.79%  "
# Row 40
$ws.Range('E40').Value = "'  +0.68%  "
# Row 41
$ws.Range('E41').Value = "'  +1.46%  "
# Row 42
$ws.Range('E42').Value = "'  -1.54%  "
# Row 43
$ws.Range('D43').Value = "'5.42"
$ws.Range('E43').Value = "'  +0.87%  "
# Row 44
$ws.Range('D44').Value = "'1.799.04"
$ws.Range('E44').Value = "'  +0.73%  "
# Row 45
$ws.Range('D45').Value = "'92.18"
$ws.Range('E45').Value = "'  -1.83%  "
# Row 46
$ws.Range('D46').Value = "'59.77"
$ws.Range('E46').Value = "'  +7.25%  "
# Row 47
$ws.Range('D47').Value = "'1.63"
$ws.Range('E47').Value = "'  +1.04%  "
# Row 48
$ws.Range('E48').Value = "'  +1.03%  "
# Row 49
$ws.Range('E49').Value = "'  +0.33%  "
# Row 50
$ws.Range('D50').Value = "'7.80"
$ws.Range('E50').Value = "'  +1.51%  "
# Row 51
$ws.Range('E51').Value = "'  +1.10%  "
